$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row at position 31, shifting existing row 31 (and below) down.
# Copy the row above first so the new row inherits the same formatting/style.
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(31).Insert()
$excel.CutCopyMode = $false

# Populate the newly inserted row 31 with the new default user.
$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"
$ws.Cells.Item(31, 4).Value = "CUSTOM_USER"
$ws.Cells.Item(31, 5).Value = "Smoke Test User"
$ws.Cells.Item(31, 6).Value = "N"

# Match the cursor/selection position seen in the authored workbook.
$ws.Range("F31").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
